$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'27.825.84"
$ws.Range("E2").Value = "'  -1.35%  "

# Row 3
$ws.Range("D3").Value = "'1.771.24"
$ws.Range("E3").Value = "'  -1.38%  "

# Row 4
$ws.Range("D4").Value = "'1.013"
$ws.Range("E4").Value = "'  +0.89%  "

# Row 5
$ws.Range("D5").Value = "'327.72"
$ws.Range("E5").Value = "'  -3.25%  "

# Row 6
$ws.Range("D6").Value = "'1.013"
$ws.Range("E6").Value = "'  +1.35%  "

# Row 7
$ws.Range("D7").Value = "'0.4325"
$ws.Range("E7").Value = "'  -5.34%  "

# Row 8
$ws.Range("D8").Value = "'0.3672"
$ws.Range("E8").Value = "'  +1.88%  "

# Row 9
$ws.Range("D9").Value = "'45.18"
$ws.Range("E9").Value = "'  -0.65%  "

# Row 10
$ws.Range("B10").Value = "'Polygon"
$ws.Range("C10").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "'1.127"
$ws.Range("E10").Value = "'  -1.08%  "

# Row 11
$ws.Range("B11").Value = "'Dogecoin"
$ws.Range("C11").Value = "'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.07451"
$ws.Range("E11").Value = "'  -0.63%  "

# Row 12
$ws.Range("D12").Value = "'1.011"
$ws.Range("E12").Value = "'  +0.86%  "

# Row 13
$ws.Range("D13").Value = "'22.01"
$ws.Range("E13").Value = "'  -1.52%  "

# Row 14
$ws.Range("D14").Value = "'6.196"
$ws.Range("E14").Value = "'  -0.25%  "

# Row 15
$ws.Range("D15").Value = "'7.353"
$ws.Range("E15").Value = "'  +1.61%  "

# Row 16
$ws.Range("D16").Value = "'1.775.65"
$ws.Range("E16").Value = "'  -0.99%  "

# Row 17
$ws.Range("D17").Value = "'0.00001075"
$ws.Range("E17").Value = "'  -0.61%  "

# Row 18
$ws.Range("D18").Value = "'0.06620"
$ws.Range("E18").Value = "'  -1.09%  "

# Row 19
$ws.Range("D19").Value = "'82.22"
$ws.Range("E19").Value = "'  +1.31%  "

# Row 20
$ws.Range("D20").Value = "'1.009"
$ws.Range("E20").Value = "'  +0.94%  "

# Row 21
$ws.Range("E21").Value = "'  -0.64%  "

# Row 22
$ws.Range("D22").Value = "'6.181"
$ws.Range("E22").Value = "'  -2.75%  "

# Row 23
$ws.Range("D23").Value = "'27.825.49"
$ws.Range("E23").Value = "'  -1.40%  "

# Row 24
$ws.Range("D24").Value = "'11.36"
$ws.Range("E24").Value = "'  -4.24%  "

# Row 25
$ws.Range("D25").Value = "'2.419"
$ws.Range("E25").Value = "'  +2.02%  "

# Row 26
$ws.Range("D26").Value = "'20.27"
$ws.Range("E26").Value = "'  -0.31%  "

# Row 27
$ws.Range("D27").Value = "'150.86"
$ws.Range("E27").Value = "'  -1.73%  "

# Row 28
$ws.Range("D28").Value = "'2.334"
$ws.Range("E28").Value = "'  -2.12%  "

# Row 29
$ws.Range("D29").Value = "'1.984.08"
$ws.Range("E29").Value = "'  -0.67%  "

# Row 30
$ws.Range("D30").Value = "'1.294"
$ws.Range("E30").Value = "'  +2.40%  "

# Row 31
$ws.Range("D31").Value = "'128.84"
$ws.Range("E31").Value = "'  -2.49%  "

# Row 32
$ws.Range("D32").Value = "'3.993"
$ws.Range("E32").Value = "'  -1.93%  "

# Row 33
$ws.Range("D33").Value = "'5.742"
$ws.Range("E33").Value = "'  -2.12%  "

# Row 34
$ws.Range("D34").Value = "'0.09131"
$ws.Range("E34").Value = "'  -3.37%  "

# Row 35
$ws.Range("E35").Value = "'  +2.65%  "

# Row 36
$ws.Range("D36").Value = "'12.08"
$ws.Range("E36").Value = "'  +0.29%  "

# Row 37
$ws.Range("D37").Value = "'0.6566"
$ws.Range("E37").Value = "'  -0.66%  "

# Row 38
$ws.Range("D38").Value = "'0.06198"
$ws.Range("E38").Value = "'  -1.02%  "

# Row 39
$ws.Range("B39").Value = "'InternetComputer(DFINITY)"
$ws.Range("C39").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'5.143"
$ws.Range("E39").Value = "'  -0.34%  "

# Row 40
$ws.Range("B40").Value = "'VeChain"
$ws.Range("C40").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.02291"
$ws.Range("E40").Value = "'  -3.02%  "

# Row 41
$ws.Range("D41").Value = "'1.198"
$ws.Range("E41").Value = "'  -1.13%  "

# Row 42
$ws.Range("D42").Value = "'1.441"
$ws.Range("E42").Value = "'  -2.66%  "

# Row 43
$ws.Range("B43").Value = "'Frax"
$ws.Range("C43").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").Value = "'1.012"
$ws.Range("E43").Value = "'  +1.26%  "

# Row 44
$ws.Range("B44").Value = "'FraxShare"
$ws.Range("C44").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'8.017"
$ws.Range("E44").Value = "'  -0.33%  "

# Row 45
$ws.Range("D45").Value = "'13.83"
$ws.Range("E45").Value = "'  -1.37%  "

# Row 46
$ws.Range("B46").Value = "'PancakeSwap"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "'3.808"
$ws.Range("E46").Value = "'  -1.56%  "

# Row 47
$ws.Range("B47").Value = "'Decentraland"
$ws.Range("C47").Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.5986"
$ws.Range("E47").Value = "'  -1.22%  "

# Row 48
$ws.Range("D48").Value = "'126.12"
$ws.Range("E48").Value = "'  -1.61%  "

# Row 49
$ws.Range("D49").Value = "'1.987"
$ws.Range("E49").Value = "'  -1.55%  "

# Row 50
$ws.Range("D50").Value = "'0.06932"
$ws.Range("E50").Value = "'  -2.26%  "

# Row 51
$ws.Range("D51").Value = "'1.123"
$ws.Range("E51").Value = "'  -3.43%  "
